$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column J: J4 = I4*2, and a shared formula J5:J7 = I5*2 / I6*2 / I7*2
$ws.Range("J4").Formula = "=I4*2"
$ws.Range("J5:J7").Formula = "=I5*2"

# Clean up row 2 styling so it matches rows 3-14 (drop the redundant
# "applyBorder" flag that was left over on these cells - border itself
# was never visible since borderId=0).
$row2Cells = @("A2", "B2", "D2", "F2")
foreach ($addr in $row2Cells) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Times New Roman"
    $c.WrapText = $true
    $c.VerticalAlignment = -4108
}

$h2 = $ws.Range("H2")
$h2.HorizontalAlignment = -4152
$h2.Font.Name = "Times New Roman"
$h2.WrapText = $true
$h2.VerticalAlignment = -4108

# Move the active selection to I4
[void]$ws.Range("I4").Select()
